$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing statuses
$ws.Range("E3").Value = "3rd Interview"
$ws.Range("E5").Value = "1st Interview"

# Insert a new row before row 8 (Laurel / Customer Success Manager UK), shifting
# the old row 8 (Accel Data / Daniel Wing) down to row 9
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 853
$ws.Range("B8").Value = "Laurel"
$ws.Range("C8").Value = "Customer Success Manager UK"
$ws.Range("D8").Value = "Regitze Steffensen"
$ws.Range("E8").Value = "CV Sent"

# Update the status for the row that shifted down to row 9 (Daniel Wing)
$ws.Range("E9").Value = "1st Interview"

# Add a new row 10 for the second Accel Data candidate
$ws.Range("A10").Value = 865
$ws.Range("B10").Value = "Accel Data"
$ws.Range("C10").Value = "Enterprise Account Executive"
$ws.Range("D10").Value = "Nicholas Lomas"
$ws.Range("E10").Value = "CV Sent"
